# Updated symbol list on Mon Dec 19 18:38:19 UTC 2022 with GitHub Actions
#
# Refreshes the "cryptos" price sheet: most rows just get a new Price
# (column D) pulled from the latest coinranking.com snapshot. Row 18 and
# row 47's rank/name label (column E) also pick up/lose a "Worst in 24h"
# suffix, and rows 42/43 swap places (BKEXToken <-> CEJI) with fresh data.
#
# Price cells store numeric-looking text (e.g. "243.86"), so we force the
# cell to Text format ("@") before assigning the new value - otherwise
# Excel would silently coerce the literal into a real number and change
# the cell's stored type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($a1, $value) {
    $rng = $ws.Range($a1)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Price (column D) refresh for most rows ---
Set-TextValue "D2"  '243.86'
Set-TextValue "D3"  '21.51'
Set-TextValue "D4"  '5.312'
Set-TextValue "D5"  '0.05637'
Set-TextValue "D6"  '3.371'
Set-TextValue "D7"  '6.380'
Set-TextValue "D8"  '0.8060'
Set-TextValue "D9"  '0.9549'
Set-TextValue "D10" '0.1426'
Set-TextValue "D11" '0.07403'
Set-TextValue "D12" '0.03212'
Set-TextValue "D13" '0.03082'
Set-TextValue "D14" '0.09276'
Set-TextValue "D15" '3.567'
Set-TextValue "D16" '0.001641'
Set-TextValue "D17" '0.04693'
Set-TextValue "D18" '0.0005811'

# Row 18 (One/ONE) label now flags "Worst in 24h"
$ws.Range("E18").Value = '17OneONEWorstin24h'

Set-TextValue "D19" '0.006348'
Set-TextValue "D20" '0.004987'
Set-TextValue "D21" '0.001043'
Set-TextValue "D23" '0.0003101'
Set-TextValue "D24" '3.767'
Set-TextValue "D25" '2.096'
Set-TextValue "D26" '0.3255'
Set-TextValue "D40" '0.03952'
Set-TextValue "D41" '0.006975'

# --- Rows 42/43 swap: BKEXToken <-> CEJI, each with refreshed data ---
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue "D42" '0.003501'
$ws.Range("E42").Value = '41CEJICEJI'

$ws.Range("B43").Value = 'BKEXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue "D43" '0.1036'
$ws.Range("E43").Value = '42BKEXTokenBKK'

Set-TextValue "D44" '0.007474'
Set-TextValue "D45" '0.00005929'
Set-TextValue "D46" '0.00000000750'
Set-TextValue "D47" '0.0005501'

# Row 47 (ACDXExchange) label loses its "Worst in 24h" flag
$ws.Range("E47").Value = '46ACDXExchangeACXT'

Set-TextValue "D48" '0.6826'
Set-TextValue "D49" '0.05638'
Set-TextValue "D50" '0.00002100'
Set-TextValue "D51" '0.01010'
